$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- 1) "Frameworks & Libraries:" -> "Frameworks:" (keep as a single run) ---
$full = $tr.Text
$needle1 = "Frameworks & Libraries:"
$idx1 = $full.IndexOf($needle1)
if ($idx1 -ge 0) {
    $sub1 = $tr.Characters($idx1 + 1, $needle1.Length)
    $sub1.Text = "Frameworks:"
}

# --- 2) "Selenium " -> "Selenium webdriver" (append new run "webdriver") ---
$full = $tr.Text
$needle2 = "Selenium "
$idx2 = $full.IndexOf($needle2)
if ($idx2 -ge 0) {
    $sub2 = $tr.Characters($idx2 + 1, $needle2.Length)
    $sub2.Text = "Selenium webdriver"
}

# Force the appended "webdriver" text into its own run (distinct a:r) while
# keeping the inherited run formatting (Georgia / sz 2000) intact, by
# touching a formatting property on just that sub-range.
$full = $tr.Text
$needle3 = "webdriver"
$idx3 = $full.IndexOf($needle3)
if ($idx3 -ge 0) {
    $sub3 = $tr.Characters($idx3 + 1, $needle3.Length)
    $sub3.Font.Size = $sub3.Font.Size
}
